$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36, pushing the existing rows 36-37 down to 37-38.
$ws.Rows("36:36").Insert()

# Copy formatting (number format / alignment / style) from neighboring row
# so the new row matches the workbook's existing look & feel.
$ws.Range("A37:C37").Copy()
$ws.Range("A36:C36").PasteSpecial(-4122)
$ws.Range("E37:F37").Copy()
$ws.Range("E36:F36").PasteSpecial(-4122)

# Row 34 is another "prospect" row (no invoice date yet) - reuse its D-column
# formatting (General number format) for the new row's blank date cell.
$ws.Range("D34").Copy()
$ws.Range("D36").PasteSpecial(-4122)

# Populate the new prospect row.
$ws.Range("A36").Value() = "MERWIN LIQUORS FALCON HEIGHTS"
$ws.Range("B36").Value() = "Larsen, Rick J"
$ws.Range("C36").Value() = "023"
$ws.Range("E36").Value() = "0008384"

# Match the custom row height used by the rest of the data rows.
$ws.Rows("36:36").RowHeight() = $ws.Rows("37:37").RowHeight()
